$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 524, pushing the existing rows 524:549 down to 525:550
$ws.Rows.Item(524).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A524").Value = 10
$ws.Range("B524").Value = "Vega Modelo de Temuco"
$ws.Range("C524").Value = "La Araucanía"
$ws.Range("D524").Value = 45041
$ws.Range("E524").Value = 9
$ws.Range("F524").Value = 100112040
$ws.Range("G524").Value = "Cilantro"
$ws.Range("H524").Value = "Sin especificar"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 50
$ws.Range("K524").Value = 5000
$ws.Range("L524").Value = 5000
$ws.Range("M524").Value = 5000
$ws.Range("N524").Value = "$/docena de atados (2 kilos)"
$ws.Range("O524").Value = "Provincia de Cautín"
$ws.Range("P524").Value = 2500
$ws.Range("Q524").Value = 2
$ws.Range("R524").Value = "Hortaliza"
